# Update results after reparse/fix builder
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newB = @(0.0965999999999999, 0.105, 0.0886, 0.1094, 0.128, 0.162, 0.2242, 0.1068, 0.3994, 0.923)
$newD = @(0.0242, 0.0322, 0.0283999999999999, 0.026, 0.0283999999999999, 0.0221999999999999, 0.0252, 0.0254, 0.0327999999999999, 0.0324)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $newB[$i]
    $ws.Cells.Item($row, 4).Value = $newD[$i]
}

[void]$ws.Range("A1:D10").Select()
